$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 08:05"

# --- Israel (row 41): new case counts ---
$ws.Range("B41").Value = 16690
$ws.Range("C41").Value = 7
$ws.Range("D41").Value = 13915
$ws.Range("E41").Value = 2496

# --- Uzbekistan (row 77): new case counts ---
$ws.Range("B77").Value = 3006
$ws.Range("C77").Value = 42
$ws.Range("E77").Value = 586

# --- Rows 81-83 re-sorted: Bulgaria moves above Tayikistan / Bosnia y Herzegovina,
#     each country keeps moving down a slot and Bulgaria gets fresh data ---

# Row 81 now holds Bulgaria's (updated) data
$ws.Range("A81").Value = "Bulgaria"
$ws.Range("B81").Value = 2372
$ws.Range("C81").Value = 41
$ws.Range("D81").Value = 769
$ws.Range("E81").Value = 1478
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 5
$ws.Range("H81").Value = 125

# Row 82 now holds Tayikistan's (previous) data
$ws.Range("A82").Value = "Tayikistan"
$ws.Range("B82").Value = 2350
$ws.Range("C82").Value = 0
$ws.Range("D82").Value = 1008
$ws.Range("E82").Value = 1298
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 44

# Row 83 now holds Bosnia y Herzegovina's (previous) data
$ws.Range("A83").Value = "Bosnia y Herzegovina"
$ws.Range("B83").Value = 2350
$ws.Range("C83").Value = 0
$ws.Range("D83").Value = 1596
$ws.Range("E83").Value = 614
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 140

# --- Rows 209-211 re-sorted: Montserrat moves above Groenlandia / Seychelles ---

# Row 209 now holds Montserrat's (previous) data
$ws.Range("A209").Value = "Montserrat"
$ws.Range("B209").Value = 11
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 10
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 1

# Row 210 now holds Groenlandia's (previous) data
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("B210").Value = 11
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 11
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Row 211 now holds Seychelles' (previous) data
$ws.Range("A211").Value = "Seychelles"
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 11
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0

# --- Rows 214-215 re-sorted: Sahara Occidental moves above Bonaire, San Eustaquio y Saba ---

# Row 214 now holds Sahara Occidental's data (values identical, only name swaps)
$ws.Range("A214").Value = "Sahara Occidental"
$ws.Range("B214").Value = 6
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 6
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

# Row 215 now holds Bonaire, San Eustaquio y Saba's data (values identical, only name swaps)
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("B215").Value = 6
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 6
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0
